$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B33").Value = "Overall Resource DCA - Now"
$ws.Range("C33").Value = "Count"
$ws.Range("D33").Value = "Costs"
$ws.Range("E33").Value = "Proportion costs"

$ws.Range("B34").Value = "Green"
$ws.Range("C34").Value = 4
$ws.Range("D34").Value = 130
$ws.Range("E34").Value = 0.8666666666666667

$ws.Range("B35").Value = "Amber/Green"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0

$ws.Range("B36").Value = "Amber"
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 20
$ws.Range("E36").Value = 0.1333333333333333

$ws.Range("B37").Value = "Amber/Red"
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0

$ws.Range("B38").Value = "Red"
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0

$ws.Range("B39").Value = "None"
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0

$ws.Range("B40").Value = "Total"
$ws.Range("C40").Value = 5
$ws.Range("D40").Value = 150
$ws.Range("E40").Value = 1
